# Scheduled-runner data refresh: updates cached market-price / profit
# figures (columns H-N) on several rows across the ALC/ARM/BSM/CRP/CUL/
# GSM/LTW/WVR sheets. A few rows also lose their trailing "HQ profit"
# cell (column M or N) because that recipe no longer has an HQ variant
# in the refreshed data pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1435.4445
$ws.Range("I2").Value = 1790.75
$ws.Range("K2").Value = 1790.75
$ws.Range("M2").Value = -1677.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 111.28571
$ws.Range("I5").Value = 46.5
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 46.5
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = 68.5
$ws.Range("N5").Value = -730

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5353.684
$ws.Range("J32").Value = 4984.4443
$ws.Range("L32").Value = 4984.4443
$ws.Range("N32").Value = -5636.4443

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 46.125
$ws.Range("I38").Value = 46.125
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 138.375
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 233.625
$ws.Range("N38").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 159.5
$ws.Range("I39").Value = 63.714287
$ws.Range("J39").Value = 293.6
$ws.Range("K39").Value = 191.142861
$ws.Range("L39").Value = 880.8000000000001
$ws.Range("M39").Value = 104.857139
$ws.Range("N39").Value = -1472.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 4985.773
$ws.Range("J53").Value = 11908.223
$ws.Range("L53").Value = 11908.223
$ws.Range("N53").Value = -13182.223

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 20833954
$ws.Range("I107").Value = 23810152
$ws.Range("J107").Value = 564.5
$ws.Range("K107").Value = 23810152
$ws.Range("L107").Value = 564.5
$ws.Range("M107").Value = -23808232
$ws.Range("N107").Value = -4404.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 35036.332
$ws.Range("J112").Value = 51554.5
$ws.Range("L112").Value = 154663.5
$ws.Range("N112").Value = -156879.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 8886.556
$ws.Range("I113").Value = 12733
$ws.Range("K113").Value = 12733
$ws.Range("M113").Value = -9479

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 5533.1816
$ws.Range("I116").Value = 4848.3335
$ws.Range("J116").Value = 6355
$ws.Range("K116").Value = 4848.3335
$ws.Range("L116").Value = 6355
$ws.Range("M116").Value = -1406.3335
$ws.Range("N116").Value = -13239

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 41670850
$ws.Range("I132").Value = 58828384
$ws.Range("K132").Value = 176485152
$ws.Range("M132").Value = -176482622

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2072.3333
$ws.Range("I141").Value = 1887.3077
$ws.Range("K141").Value = 5661.9231
$ws.Range("M141").Value = -481.9231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2925390.8
$ws.Range("I2").Value = 5051549.5
$ws.Range("J2").Value = 1922.125
$ws.Range("K2").Value = 5051549.5
$ws.Range("L2").Value = 1922.125
$ws.Range("M2").Value = -5051436.5
$ws.Range("N2").Value = -2148.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 10279493
$ws.Range("I45").Value = 23977372
$ws.Range("J45").Value = 6083.375
$ws.Range("K45").Value = 23977372
$ws.Range("L45").Value = 6083.375
$ws.Range("M45").Value = -23976995
$ws.Range("N45").Value = -6837.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2861.182
$ws.Range("I61").Value = 2481
$ws.Range("J61").Value = 3875
$ws.Range("K61").Value = 2481
$ws.Range("L61").Value = 3875
$ws.Range("M61").Value = -2269
$ws.Range("N61").Value = -4299

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 52045.918
$ws.Range("I74").Value = 2716.4062
$ws.Range("K74").Value = 2716.4062
$ws.Range("M74").Value = -1842.4062

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 52045.918
$ws.Range("I77").Value = 2716.4062
$ws.Range("K77").Value = 13582.031
$ws.Range("M77").Value = -9214.030999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2925390.8
$ws.Range("I116").Value = 5051549.5
$ws.Range("J116").Value = 1922.125
$ws.Range("K116").Value = 5051549.5
$ws.Range("L116").Value = 1922.125
$ws.Range("M116").Value = -5049255.5
$ws.Range("N116").Value = -6510.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1911.4849
$ws.Range("I132").Value = 1607.4
$ws.Range("J132").Value = 2861.75
$ws.Range("K132").Value = 4822.200000000001
$ws.Range("L132").Value = 8585.25
$ws.Range("M132").Value = -2292.200000000001
$ws.Range("N132").Value = -13645.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2861.182
$ws.Range("I136").Value = 2481
$ws.Range("J136").Value = 3875
$ws.Range("K136").Value = 7443
$ws.Range("L136").Value = 11625
$ws.Range("M136").Value = -4893
$ws.Range("N136").Value = -16725

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2925390.8
$ws.Range("I3").Value = 5051549.5
$ws.Range("J3").Value = 1922.125
$ws.Range("K3").Value = 5051549.5
$ws.Range("L3").Value = 1922.125
$ws.Range("M3").Value = -5051435.5
$ws.Range("N3").Value = -2150.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2717.9614
$ws.Range("I134").Value = 1388.1163
$ws.Range("J134").Value = 9071.666999999999
$ws.Range("K134").Value = 4164.3489
$ws.Range("L134").Value = 27215.001
$ws.Range("M134").Value = -1629.3489
$ws.Range("N134").Value = -32285.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14267.862
$ws.Range("I31").Value = 2088.182
$ws.Range("J31").Value = 18390.215
$ws.Range("K31").Value = 2088.182
$ws.Range("L31").Value = 18390.215
$ws.Range("M31").Value = -1793.182
$ws.Range("N31").Value = -18980.215

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 14267.862
$ws.Range("I34").Value = 2088.182
$ws.Range("J34").Value = 18390.215
$ws.Range("K34").Value = 2088.182
$ws.Range("L34").Value = 18390.215
$ws.Range("M34").Value = -1886.182
$ws.Range("N34").Value = -18794.215

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4090
$ws.Range("I99").Value = 3750
$ws.Range("J99").Value = 4284.2856
$ws.Range("K99").Value = 3750
$ws.Range("L99").Value = 4284.2856
$ws.Range("M99").Value = -2252
$ws.Range("N99").Value = -7280.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4090
$ws.Range("I126").Value = 3750
$ws.Range("J126").Value = 4284.2856
$ws.Range("K126").Value = 11250
$ws.Range("L126").Value = 12852.8568
$ws.Range("M126").Value = -8780
$ws.Range("N126").Value = -17792.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3308.125
$ws.Range("I134").Value = 2115
$ws.Range("J134").Value = 4842.143
$ws.Range("K134").Value = 6345
$ws.Range("L134").Value = 14526.429
$ws.Range("M134").Value = -3810
$ws.Range("N134").Value = -19596.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 811.0909
$ws.Range("I3").Value = 811.0909
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2433.2727
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -2321.2727
$ws.Range("N3").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 396.85715
$ws.Range("J10").Value = 2000
$ws.Range("L10").Value = 6000
$ws.Range("N10").Value = -6278

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 68582.46000000001
$ws.Range("I12").Value = 98891.78
$ws.Range("J12").Value = 386.5
$ws.Range("K12").Value = 296675.34
$ws.Range("L12").Value = 1159.5
$ws.Range("M12").Value = -296502.34
$ws.Range("N12").Value = -1505.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 1341.3125
$ws.Range("I108").Value = 1589
$ws.Range("K108").Value = 4767
$ws.Range("M108").Value = -1887

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3290.2917
$ws.Range("J113").Value = 2248.35
$ws.Range("L113").Value = 6745.049999999999
$ws.Range("N113").Value = -11085.05

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2222877.5
$ws.Range("I129").Value = 2500562.2
$ws.Range("J129").Value = 1400
$ws.Range("K129").Value = 7501686.600000001
$ws.Range("L129").Value = 4200
$ws.Range("M129").Value = -7496686.600000001
$ws.Range("N129").Value = -14200

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3270.1333
$ws.Range("I134").Value = 1465.5385
$ws.Range("K134").Value = 4396.6155
$ws.Range("M134").Value = 673.3845000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 1554.1111
$ws.Range("I136").Value = 1554.1111
$ws.Range("K136").Value = 4662.3333
$ws.Range("M136").Value = 437.6666999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1913.4
$ws.Range("I139").Value = 1266.75
$ws.Range("K139").Value = 3800.25
$ws.Range("M139").Value = 1339.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1018789.7
$ws.Range("I80").Value = 1743688.4
$ws.Range("K80").Value = 1743688.4
$ws.Range("M80").Value = -1742690.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 1018789.7
$ws.Range("I83").Value = 1743688.4
$ws.Range("K83").Value = 8718442
$ws.Range("M83").Value = -8713450

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1243.8334
$ws.Range("I107").Value = 1266.6666
$ws.Range("J107").Value = 1198.1666
$ws.Range("K107").Value = 1266.6666
$ws.Range("L107").Value = 1198.1666
$ws.Range("M107").Value = 653.3334
$ws.Range("N107").Value = -5038.1666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I113").Value = 83335840
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 83335840
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -83333670
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 90468.7
$ws.Range("I22").Value = 888888
$ws.Range("K22").Value = 888888
$ws.Range("M22").Value = -888593

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 90468.7
$ws.Range("I27").Value = 888888
$ws.Range("K27").Value = 888888
$ws.Range("M27").Value = -888781

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6846.4
$ws.Range("I46").Value = 4719.9
$ws.Range("J46").Value = 11099.4
$ws.Range("K46").Value = 4719.9
$ws.Range("L46").Value = 11099.4
$ws.Range("M46").Value = -4531.9
$ws.Range("N46").Value = -11475.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 11114390
$ws.Range("I61").Value = 11114390
$ws.Range("K61").Value = 11114390
$ws.Range("M61").Value = -11114188

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 11114390
$ws.Range("I113").Value = 11114390
$ws.Range("K113").Value = 11114390
$ws.Range("M113").Value = -11112220

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 43338.02
$ws.Range("I136").Value = 56224.406
$ws.Range("K136").Value = 168673.218
$ws.Range("M136").Value = -166123.218

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1322.5
$ws.Range("I113").Value = 872
$ws.Range("J113").Value = 1998.25
$ws.Range("K113").Value = 2616
$ws.Range("L113").Value = 5994.75
$ws.Range("M113").Value = -446
$ws.Range("N113").Value = -10334.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 28887058
$ws.Range("I132").Value = 37042410
$ws.Range("K132").Value = 111127230
$ws.Range("M132").Value = -111124700
